$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "The Game Studies Practicum: Applying Situated Learning to Teach Professional Practices"
$ws.Range("B20").Value = "Käsittelee Singapore-Mit GAMBIT kesäohjelmaa, miten sitä on järjestetty ja pyritty opettamaan. Näkemystä Situated learningista ja CA:sta. 3/5"
$ws.Range("C20").Value = "Singapore-MIT GAMBIT kesäohjelma 8 viikon mittainen ohjelmointiprojekti oppilailla (ohtuprojekti fiilis). Mutta opettamassa ohjelmointialan ammattilaisia pelientekonäkökulmasta. Ohjaavat oppilaita oikeaan suuntaan coachin. Kertoo situated learningista ja CA:sta. Selittää miten kesäohjelmaa on järjestetty ja miten siinä näkyvät nämä eri näkökulmat."
$ws.Range("D20").Value = "Vähän erilaista näkemystä CA:han, koska yhdistetty scrum ja ei oikein muuta CA yhteyttä. Scrumiin yhdistettyä CA ja kuinka pyritty opettamaan tätä Singapore-MIT GAMBIT kurssia. Samalla tosin epäilyttävää."
$ws.Range("E20").Value = "Ei tuloksia, kerrotaan vain miten kurssia on järjestetty ja vähän omituinen yhteys scrum liitteestä CA:han."
$ws.Range("F20").Value = "Singapore, yliopisto, (yhteistyössä MIT)"
$ws.Range("G20").Value = "Ei tuloksia"
